# Updates the Price (D) and Volume(1h) (E) columns of the cryptos table
# with freshly scraped values (GitHub Actions cron refresh).
#
# Both columns hold plain text in the workbook (e.g. "568.55", "3.396.26"
# multi-dot big numbers, "  -1.40%  " padded percentages) rather than real
# numbers/percentages. Assigning a bare numeric-looking string via .Value
# would make Excel auto-convert it to a Double, so those cells are written
# through .Formula with a leading apostrophe to force a text literal - the
# same trick Excel itself uses when a user types a number as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.121.90"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "3.387.86"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Formula = "'568.55"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Formula = "'155.39"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("D7").Formula = "'0.625"
$ws.Range("E7").Value = "  +8.06%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.390.11"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").Formula = "'7.11"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Formula = "'0.436"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "3.976.38"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Formula = "'0.0000186"
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").Formula = "'27.23"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "64.176.39"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "3.378.42"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").Formula = "'6.25"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").Formula = "'13.74"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("D21").Formula = "'375.36"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Formula = "'7.93"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Formula = "'0.539"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Formula = "'71.46"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("E26").Value = "  -5.56%  "
$ws.Range("D27").Formula = "'10.28"
$ws.Range("E27").Value = "  +6.43%  "
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Formula = "'1.45"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Formula = "'6.07"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("D33").Formula = "'22.96"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").Formula = "'7.05"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +5.85%  "
$ws.Range("D36").Formula = "'160.27"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").Formula = "'0.0753"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").Value = "2.857.32"
$ws.Range("E39").Value = "  -6.02%  "
$ws.Range("D40").Formula = "'6.72"
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").Formula = "'26.18"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").Formula = "'4.55"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Formula = "'42.57"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Formula = "'0.0311"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").Formula = "'0.764"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Formula = "'25.65"
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("D47").Formula = "'319.56"
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Formula = "'1.05"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Formula = "'6.48"
$ws.Range("E51").Value = "  -2.01%  "
